# Scheduled runner update: refresh currentAveragePrice/Price/Profit columns
# (H:N) for a handful of leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR
# sheets, reflecting newer market-board pulls.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 66.5
$ws.Range("I9").Value = 66.5
$ws.Range("K9").Value = 66.5
$ws.Range("M9").Value = 102.5
$ws.Range("H13").Value = 602.5
$ws.Range("I13").Value = 5
$ws.Range("J13").Value = 1200
$ws.Range("K13").Value = 5
$ws.Range("L13").Value = 1200
$ws.Range("M13").Value = 164
$ws.Range("N13").Value = -1538
$ws.Range("H70").Value = 1766.6666
$ws.Range("I70").Value = 500
$ws.Range("J70").Value = 2400
$ws.Range("K70").Value = 1500
$ws.Range("L70").Value = 7200
$ws.Range("M70").Value = -1230
$ws.Range("N70").Value = -7740
$ws.Range("H73").Value = 1766.6666
$ws.Range("I73").Value = 500
$ws.Range("J73").Value = 2400
$ws.Range("K73").Value = 1500
$ws.Range("L73").Value = 7200
$ws.Range("M73").Value = -564
$ws.Range("N73").Value = -9072
$ws.Range("H111").Value = 4524.143
$ws.Range("I111").Value = 6669
$ws.Range("J111").Value = 4166.6665
$ws.Range("K111").Value = 20007
$ws.Range("L111").Value = 12499.9995
$ws.Range("M111").Value = -16940
$ws.Range("N111").Value = -18633.9995
$ws.Range("H137").Value = 914641.4399999999
$ws.Range("I137").Value = 1114117.4
$ws.Range("J137").Value = 17000
$ws.Range("K137").Value = 3342352.2
$ws.Range("L137").Value = 51000
$ws.Range("M137").Value = -3339802.2
$ws.Range("N137").Value = -56100

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 163.55556
$ws.Range("I4").Value = 171.625
$ws.Range("K4").Value = 171.625
$ws.Range("M4").Value = -55.625
$ws.Range("H5").Value = 241
$ws.Range("I5").Value = 241
$ws.Range("K5").Value = 241
$ws.Range("M5").Value = -129
$ws.Range("H17").Value = 5000
$ws.Range("J17").Value = 5000
$ws.Range("L17").Value = 5000
$ws.Range("N17").Value = -5346
$ws.Range("H63").Value = 6689.6
$ws.Range("J63").Value = 3999
$ws.Range("L63").Value = 3999
$ws.Range("N63").Value = -5371
$ws.Range("H66").Value = 6689.6
$ws.Range("J66").Value = 3999
$ws.Range("L66").Value = 19995
$ws.Range("N66").Value = -26859
$ws.Range("H97").Value = 1265.6666
$ws.Range("I97").Value = 1265.6666
$ws.Range("K97").Value = 1265.6666
$ws.Range("M97").Value = -769.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 241
$ws.Range("I4").Value = 241
$ws.Range("K4").Value = 241
$ws.Range("M4").Value = -126
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 87.27273
$ws.Range("I7").Value = 38.75
$ws.Range("J7").Value = 216.66667
$ws.Range("K7").Value = 38.75
$ws.Range("L7").Value = 216.66667
$ws.Range("M7").Value = 74.25
$ws.Range("N7").Value = -442.66667
$ws.Range("H16").Value = 858.5
$ws.Range("I16").Value = 858.5
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 858.5
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -571.5
$ws.Range("N16").ClearContents()
$ws.Range("H22").Value = 311.58334
$ws.Range("I22").Value = 292.375
$ws.Range("K22").Value = 292.375
$ws.Range("M22").Value = 57.625
$ws.Range("H88").Value = 14779.429
$ws.Range("J88").Value = 14779.429
$ws.Range("L88").Value = 14779.429
$ws.Range("N88").Value = -15591.429
$ws.Range("H91").Value = 14779.429
$ws.Range("J91").Value = 14779.429
$ws.Range("L91").Value = 14779.429
$ws.Range("N91").Value = -17587.429
$ws.Range("H113").Value = 858.5
$ws.Range("I113").Value = 858.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 858.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1311.5
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("M109").ClearContents()
$ws.Range("I122").Value = 542
$ws.Range("J122").Value = 382.33334
$ws.Range("K122").Value = 4878
$ws.Range("L122").Value = 3441.00006
$ws.Range("M122").Value = -2428
$ws.Range("N122").Value = -8341.00006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 67.375
$ws.Range("I2").Value = 67.375
$ws.Range("K2").Value = 67.375
$ws.Range("M2").Value = 45.625
$ws.Range("H94").Value = 50000
$ws.Range("J94").Value = 50000
$ws.Range("L94").Value = 50000
$ws.Range("N94").Value = -51352

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8000
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 8000
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 8000
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -8224
$ws.Range("H68").Value = 4000
$ws.Range("J68").Value = 5000
$ws.Range("L68").Value = 5000
$ws.Range("N68").Value = -6498
$ws.Range("H71").Value = 4000
$ws.Range("J71").Value = 5000
$ws.Range("L71").Value = 25000
$ws.Range("N71").Value = -32488
$ws.Range("H93").Value = 1939.4
$ws.Range("I93").Value = 1674.5
$ws.Range("K93").Value = 1674.5
$ws.Range("M93").Value = -426.5
$ws.Range("H100").Value = 1618.5
$ws.Range("I100").Value = 1618.5
$ws.Range("K100").Value = 1618.5
$ws.Range("M100").Value = -1077.5
$ws.Range("H105").Value = 35000
$ws.Range("J105").Value = 35000
$ws.Range("L105").Value = 35000
$ws.Range("N105").Value = -41988
$ws.Range("H122").Value = 3889.8
$ws.Range("I122").Value = 3299.3333
$ws.Range("J122").Value = 4142.857
$ws.Range("K122").Value = 9897.999899999999
$ws.Range("L122").Value = 12428.571
$ws.Range("M122").Value = -7447.999899999999
$ws.Range("N122").Value = -17328.571
$ws.Range("H126").Value = 8000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 8000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 24000
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -28940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("H126").Value = 5000
$ws.Range("J126").Value = 5000
$ws.Range("L126").Value = 15000
$ws.Range("N126").Value = -19940
